$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fasta-method-1")
# fasta-method-1 is also the workbook's active sheet, but we address it by
# name explicitly so the script is robust regardless of which sheet/tab
# happens to be active when it runs.

$ws.Range("C2").Value = 90.856513023399998
$ws.Range("D2").Value = 90.872089862799996
$ws.Range("E2").Value = 38.98828125
$ws.Range("F2").Value = 168.953125

$ws.Range("C3").Value = 18.001127004600001
$ws.Range("D3").Value = 108.873216867
$ws.Range("E3").Value = 49.51171875
$ws.Range("F3").Value = 179.46875

$ws.Range("C4").Value = 0.00136613845825
$ws.Range("D4").Value = 49258.458025
$ws.Range("E4").Value = 102.07421875
$ws.Range("F4").Value = 233.19921875

$ws.Range("C5").Value = 0.0037100315094
$ws.Range("D5").Value = 62391.398842800001
$ws.Range("E5").Value = 102.125
$ws.Range("F5").Value = 233.15234375

$ws.Activate()
$ws.Range("I7").Select()
